# Autogenerated on Sun Feb 01 2015 22:24:41 GMT-0500 (Eastern Standard Time)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Data" to "Summary"
$ws.Name = "Summary"

# Remove the old header/label row (content moves further down the sheet)
$ws.Range("B5:D5").Clear()
$ws.Range("A6").Clear()

# New source-type note under the title, italic (matches existing "title_" / source style)
$ws.Range("A9").Value = "Source Type: SME Associations"
$ws.Range("A9").Font.Italic = $true

# Move the MSME/SME/Micro header row down from row 5 to row 11 (bold, as before)
$ws.Range("B11").Value = "Micro"
$ws.Range("B11").Font.Bold = $true

$ws.Range("C11").Value = "SMEs"
$ws.Range("C11").Font.Bold = $true

$ws.Range("D11").Value = "MSMEs"
$ws.Range("D11").Font.Bold = $true

# Move the "Enterprises (% of total)" label down from row 6 to row 12 (bold, as before)
$ws.Range("A12").Value = "Enterprises (% of total)"
$ws.Range("A12").Font.Bold = $true

# New data value next to it, kept as plain text (not a number), Normal style
$ws.Range("D12").Value = "'90"

# New source citation row, bold + underlined (new style)
$ws.Range("A13").Value = "Source: WAMDA, 2007"
$ws.Range("A13").Font.Bold = $true
$ws.Range("A13").Font.Underline = $true
